$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 (b.md) moves from "Handed back: in sync with en-US"
# to "Ready for handoff", with a fresh generate-date timestamp.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 00:35:10"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets a new handoff file/datetime, duplicate flag
# flips to False, status becomes "Ready for handoff", and an error message is
# recorded about the stale handback file.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text (otherwise Excel auto-converts "False" to a
# Boolean); resetting the style afterwards drops the quote-prefix flag that
# trick leaves behind so the cell format matches its neighbours.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 00:35:00"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c714b3fe6f84f1b945c574661a2d33c0023732b7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2347e554869c7634a13a430e907de34303215741/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: same shape of update as zh-cn, but with the de-de handoff file
# name and its own handoff datetime.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 00:35:10"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c714b3fe6f84f1b945c574661a2d33c0023732b7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2347e554869c7634a13a430e907de34303215741/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
